# Scrape Global => consolidated_data
# The tender list was re-scraped: the header label was corrected and the
# table rows shifted (one closed tender dropped off, later rows moved up,
# and the "Entretien des espaces verts" tender reappeared lower down the
# list with a new closing date).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header rename: object -> objet
$ws.Range("A1").Value = "objet"

# Rows 5-9 shift up (old row 5 "Entretien des espaces verts" / 19/06/2025 02:30
# drops out at the top of the list, and reappears at row 9 with a new date).
$ws.Range("A5").Value = "AO short Term Lease of 01 CF34-10E6G07 engine"
$ws.Range("B5").Value = "20/06/2025 20:30"

$ws.Range("A6").Value = "AO confection housses sièges B787"
$ws.Range("B6").Value = "20/06/2025 20:30"

$ws.Range("A7").Value = "Assistance comptable"
$ws.Range("B7").Value = "23/06/2025 20:30"

$ws.Range("A8").Value = "Remplacement Firewalls en fin de vie et gestion du spare"
$ws.Range("B8").Value = "25/06/2025 20:30"

$ws.Range("A9").Value = "Entretien des espaces verts"
$ws.Range("B9").Value = "26/06/2025 02:30"
